$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2023-09-14 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-15 Friday", 2)

# Update the division problems in the table, addressed by (row, column)
# so that duplicate values (e.g. "99÷4=") are handled unambiguously.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "61÷6="
$t.Cell(1, 2).Range.Text = "41÷3="
$t.Cell(1, 3).Range.Text = "11÷7="
$t.Cell(1, 4).Range.Text = "40÷8="
$t.Cell(1, 5).Range.Text = "39÷4="

$t.Cell(5, 1).Range.Text = "41÷8="
$t.Cell(5, 2).Range.Text = "21÷6="
$t.Cell(5, 3).Range.Text = "23÷7="
$t.Cell(5, 4).Range.Text = "41÷4="
$t.Cell(5, 5).Range.Text = "60÷9="

$t.Cell(9, 1).Range.Text = "24÷3="
$t.Cell(9, 2).Range.Text = "69÷5="
$t.Cell(9, 3).Range.Text = "97÷4="
$t.Cell(9, 4).Range.Text = "66÷2="
$t.Cell(9, 5).Range.Text = "11÷6="

$t.Cell(13, 1).Range.Text = "81÷6="
$t.Cell(13, 2).Range.Text = "52÷7="
$t.Cell(13, 3).Range.Text = "58÷6="
$t.Cell(13, 4).Range.Text = "41÷3="
$t.Cell(13, 5).Range.Text = "10÷5="

$t.Cell(17, 1).Range.Text = "79÷7="
$t.Cell(17, 2).Range.Text = "56÷8="
$t.Cell(17, 3).Range.Text = "92÷6="
$t.Cell(17, 4).Range.Text = "17÷4="
$t.Cell(17, 5).Range.Text = "87÷7="
